{"js": "const replacements = [\n  [\"2024-05-16 Thursday\", \"2024-05-17 Friday\"],\n  [\"710\u00d79=\", \"531\u00d74=\"],\n  [\"723\u00d72=\", \"586\u00d74=\"],\n  [\"381\u00d75=\", \"787\u00d73=\"],\n  [\"419\u00d77=\", \"178\u00d72=\"],\n  [\"799\u00d78=\", \"947\u00d78=\"],\n  [\"226\u00d78=\", \"361\u00d79=\"],\n  [\"922\u00d74=\", \"250\u00d76=\"],\n  [\"652\u00d75=\", \"833\u00d79=\"],\n  [\"720\u00d77=\", \"802\u00d76=\"],\n  [\"764\u00d77=\", \"619\u00d76=\"],\n  [\"814\u00d79=\", \"484\u00d78=\"],\n  [\"165\u00d74=\", \"125\u00d79=\"],\n  [\"265\u00d74=\", \"321\u00d77=\"],\n  [\"231\u00d73=\", \"691\u00d72=\"],\n  [\"337\u00d79=\", \"543\u00d77=\"],\n  [\"522\u00d72=\", \"781\u00d77=\"],\n  [\"683\u00d75=\", \"509\u00d72=\"],\n  [\"893\u00d77=\", \"568\u00d73=\"],\n  [\"329\u00d79=\", \"667\u00d73=\"],\n  [\"618\u00d79=\", \"417\u00d78=\"],\n  [\"253\u00d77=\", \"766\u00d72=\"],\n  [\"779\u00d74=\", \"572\u00d79=\"],\n  [\"399\u00d79=\", \"402\u00d77=\"],\n  [\"376\u00d75=\", \"157\u00d77=\"],\n  [\"326\u00d77=\", \"632\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Mapping of old text -> new text (date line + multiplication problems)\n$replacements = @{\n    \"2024-05-16 Thursday\" = \"2024-05-17 Friday\"\n    \"710\u00d79=\"               = \"531\u00d74=\"\n    \"723\u00d72=\"               = \"586\u00d74=\"\n    \"381\u00d75=\"               = \"787\u00d73=\"\n    \"419\u00d77=\"               = \"178\u00d72=\"\n    \"799\u00d78=\"               = \"947\u00d78=\"\n    \"226\u00d78=\"               = \"361\u00d79=\"\n    \"922\u00d74=\"               = \"250\u00d76=\"\n    \"652\u00d75=\"               = \"833\u00d79=\"\n    \"720\u00d77=\"               = \"802\u00d76=\"\n    \"764\u00d77=\"               = \"619\u00d76=\"\n    \"814\u00d79=\"               = \"484\u00d78=\"\n    \"165\u00d74=\"               = \"125\u00d79=\"\n    \"265\u00d74=\"               = \"321\u00d77=\"\n    \"231\u00d73=\"               = \"691\u00d72=\"\n    \"337\u00d79=\"               = \"543\u00d77=\"\n    \"522\u00d72=\"               = \"781\u00d77=\"\n    \"683\u00d75=\"               = \"509\u00d72=\"\n    \"893\u00d77=\"               = \"568\u00d73=\"\n    \"329\u00d79=\"               = \"667\u00d73=\"\n    \"618\u00d79=\"               = \"417\u00d78=\"\n    \"253\u00d77=\"               = \"766\u00d72=\"\n    \"779\u00d74=\"               = \"572\u00d79=\"\n    \"399\u00d79=\"               = \"402\u00d77=\"\n    \"376\u00d75=\"               = \"157\u00d77=\"\n    \"326\u00d77=\"               = \"632\u00d75=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
